$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 54.69462833333333
$ws.Range("H2").Value = 164.083885
$ws.Range("I2").Value = 0.2790924419198448
$ws.Range("J2").Value = 0.2790924419198448
$ws.Range("M2").Value = 14.861848
$ws.Range("N2").Value = 44.585544
$ws.Range("O2").Value = 0.09055189482833943
$ws.Range("P2").Value = 0.09055189482833945
$ws.Range("Q2").Value = 812.8632527064933
$ws.Range("R2").Value = 7315.76927435844
$ws.Range("S2").Value = 0.02527234944811022
$ws.Range("T2").Value = 0.02527234944811022
$ws.Range("G3").Value = 54.69462833333333
$ws.Range("H3").Value = 164.083885
$ws.Range("I3").Value = 0.2790924419198448
$ws.Range("J3").Value = 0.2790924419198448
$ws.Range("O3").Value = 0.1893562842131466
$ws.Range("P3").Value = 0.1893562842131466
$ws.Range("Q3").Value = 1699.807225433586
$ws.Range("R3").Value = 15298.26502890227
$ws.Range("S3").Value = 0.05284790775391523
$ws.Range("T3").Value = 0.05284790775391524
$ws.Range("G4").Value = 54.69462833333333
$ws.Range("H4").Value = 164.083885
$ws.Range("I4").Value = 0.2790924419198448
$ws.Range("J4").Value = 0.2790924419198448
$ws.Range("M4").Value = 18.10188466666667
$ws.Range("N4").Value = 54.305654
$ws.Range("O4").Value = 0.1102931450066459
$ws.Range("P4").Value = 0.1102931450066459
$ws.Range("Q4").Value = 990.0758539761989
$ws.Range("R4").Value = 8910.68268578579
$ws.Range("S4").Value = 0.03078198316692433
$ws.Range("T4").Value = 0.03078198316692434
$ws.Range("G5").Value = 54.69462833333333
$ws.Range("H5").Value = 164.083885
$ws.Range("I5").Value = 0.2790924419198448
$ws.Range("J5").Value = 0.2790924419198448
$ws.Range("M5").Value = 100.0833306666667
$ws.Range("N5").Value = 300.249992
$ws.Range("O5").Value = 0.609798675951868
$ws.Range("P5").Value = 0.6097986759518681
$ws.Range("Q5").Value = 5474.020573175435
$ws.Range("R5").Value = 49266.18515857892
$ws.Range("S5").Value = 0.170190201550895
$ws.Range("T5").Value = 0.170190201550895
$ws.Range("G6").Value = 19.32115333333334
$ws.Range("H6").Value = 57.96346000000001
$ws.Range("I6").Value = 0.09859081282432611
$ws.Range("J6").Value = 0.09859081282432611
$ws.Range("M6").Value = 14.861848
$ws.Range("N6").Value = 44.585544
$ws.Range("O6").Value = 0.09055189482833943
$ws.Range("P6").Value = 0.09055189482833945
$ws.Range("Q6").Value = 287.1480440246934
$ws.Range("R6").Value = 2584.33239622224
$ws.Range("S6").Value = 0.008927584913908876
$ws.Range("T6").Value = 0.008927584913908878
$ws.Range("G7").Value = 19.32115333333334
$ws.Range("H7").Value = 57.96346000000001
$ws.Range("I7").Value = 0.09859081282432611
$ws.Range("J7").Value = 0.09859081282432611
$ws.Range("O7").Value = 0.1893562842131466
$ws.Range("P7").Value = 0.1893562842131466
$ws.Range("Q7").Value = 600.4654760528779
$ws.Range("R7").Value = 5404.189284475901
$ws.Range("S7").Value = 0.01866878997396823
$ws.Range("T7").Value = 0.01866878997396823
$ws.Range("G8").Value = 19.32115333333334
$ws.Range("H8").Value = 57.96346000000001
$ws.Range("I8").Value = 0.09859081282432611
$ws.Range("J8").Value = 0.09859081282432611
$ws.Range("M8").Value = 18.10188466666667
$ws.Range("N8").Value = 54.305654
$ws.Range("O8").Value = 0.1102931450066459
$ws.Range("P8").Value = 0.1102931450066459
$ws.Range("Q8").Value = 349.7492892669823
$ws.Range("R8").Value = 3147.743603402841
$ws.Range("S8").Value = 0.01087389081515648
$ws.Range("T8").Value = 0.01087389081515648
$ws.Range("G9").Value = 19.32115333333334
$ws.Range("H9").Value = 57.96346000000001
$ws.Range("I9").Value = 0.09859081282432611
$ws.Range("J9").Value = 0.09859081282432611
$ws.Range("M9").Value = 100.0833306666667
$ws.Range("N9").Value = 300.249992
$ws.Range("O9").Value = 0.609798675951868
$ws.Range("P9").Value = 0.6097986759518681
$ws.Range("Q9").Value = 1933.725377921369
$ws.Range("R9").Value = 17403.52840129232
$ws.Range("S9").Value = 0.06012054712129251
$ws.Range("T9").Value = 0.06012054712129251
$ws.Range("G10").Value = 11.023718
$ws.Range("H10").Value = 33.071154
$ws.Range("I10").Value = 0.05625116157486912
$ws.Range("J10").Value = 0.05625116157486911
$ws.Range("M10").Value = 14.861848
$ws.Range("N10").Value = 44.585544
$ws.Range("O10").Value = 0.09055189482833943
$ws.Range("P10").Value = 0.09055189482833945
$ws.Range("Q10").Value = 163.832821310864
$ws.Range("R10").Value = 1474.495391797776
$ws.Range("S10").Value = 0.005093649266899477
$ws.Range("T10").Value = 0.005093649266899477
$ws.Range("G11").Value = 11.023718
$ws.Range("H11").Value = 33.071154
$ws.Range("I11").Value = 0.05625116157486912
$ws.Range("J11").Value = 0.05625116157486911
$ws.Range("O11").Value = 0.1893562842131466
$ws.Range("P11").Value = 0.1893562842131466
$ws.Range("Q11").Value = 342.5966329516567
$ws.Range("R11").Value = 3083.36969656491
$ws.Range("S11").Value = 0.01065151093849055
$ws.Range("T11").Value = 0.01065151093849055
$ws.Range("G12").Value = 11.023718
$ws.Range("H12").Value = 33.071154
$ws.Range("I12").Value = 0.05625116157486912
$ws.Range("J12").Value = 0.05625116157486911
$ws.Range("M12").Value = 18.10188466666667
$ws.Range("N12").Value = 54.305654
$ws.Range("O12").Value = 0.1102931450066459
$ws.Range("P12").Value = 0.1102931450066459
$ws.Range("Q12").Value = 199.5500718338573
$ws.Range("R12").Value = 1795.950646504716
$ws.Range("S12").Value = 0.006204117520369307
$ws.Range("T12").Value = 0.006204117520369308
$ws.Range("G13").Value = 11.023718
$ws.Range("H13").Value = 33.071154
$ws.Range("I13").Value = 0.05625116157486912
$ws.Range("J13").Value = 0.05625116157486911
$ws.Range("M13").Value = 100.0833306666667
$ws.Range("N13").Value = 300.249992
$ws.Range("O13").Value = 0.609798675951868
$ws.Range("P13").Value = 0.6097986759518681
$ws.Range("Q13").Value = 1103.290413770085
$ws.Range("R13").Value = 9929.613723930766
$ws.Range("S13").Value = 0.03430188384910978
$ws.Range("T13").Value = 0.03430188384910978
$ws.Range("G14").Value = 110.9336623333333
$ws.Range("H14").Value = 332.800987
$ws.Range("I14").Value = 0.5660655836809599
$ws.Range("J14").Value = 0.5660655836809599
$ws.Range("M14").Value = 14.861848
$ws.Range("N14").Value = 44.585544
$ws.Range("O14").Value = 0.09055189482833943
$ws.Range("P14").Value = 0.09055189482833945
$ws.Range("Q14").Value = 1648.679227681325
$ws.Range("R14").Value = 14838.11304913193
$ws.Range("S14").Value = 0.05125831119942086
$ws.Range("T14").Value = 0.05125831119942087
$ws.Range("G15").Value = 110.9336623333333
$ws.Range("H15").Value = 332.800987
$ws.Range("I15").Value = 0.5660655836809599
$ws.Range("J15").Value = 0.5660655836809599
$ws.Range("O15").Value = 0.1893562842131466
$ws.Range("P15").Value = 0.1893562842131466
$ws.Range("Q15").Value = 3447.611703818622
$ws.Range("R15").Value = 31028.5053343676
$ws.Range("S15").Value = 0.1071880755467726
$ws.Range("T15").Value = 0.1071880755467726
$ws.Range("G16").Value = 110.9336623333333
$ws.Range("H16").Value = 332.800987
$ws.Range("I16").Value = 0.5660655836809599
$ws.Range("J16").Value = 0.5660655836809599
$ws.Range("M16").Value = 18.10188466666667
$ws.Range("N16").Value = 54.305654
$ws.Range("O16").Value = 0.1102931450066459
$ws.Range("P16").Value = 0.1102931450066459
$ws.Range("Q16").Value = 2008.108361208944
$ws.Range("R16").Value = 18072.9752508805
$ws.Range("S16").Value = 0.06243315350419576
$ws.Range("T16").Value = 0.06243315350419577
$ws.Range("G17").Value = 110.9336623333333
$ws.Range("H17").Value = 332.800987
$ws.Range("I17").Value = 0.5660655836809599
$ws.Range("J17").Value = 0.5660655836809599
$ws.Range("M17").Value = 100.0833306666667
$ws.Range("N17").Value = 300.249992
$ws.Range("O17").Value = 0.609798675951868
$ws.Range("P17").Value = 0.6097986759518681
$ws.Range("Q17").Value = 11102.61040937134
$ws.Range("R17").Value = 99923.49368434209
$ws.Range("S17").Value = 0.3451860434305707
$ws.Range("T17").Value = 0.3451860434305707
